$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold plain numeric-looking strings (e.g. "51.66")
# that Excel would otherwise auto-coerce into real numbers on assignment.
# Force them to Text format first, write the value, then restore the
# cell style so no stray formatting is introduced.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "51.987.88"
$ws.Range("E2").Value = "  +0.73%  "
Set-TextValue $ws.Range("D3") "2.790.29"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue $ws.Range("D5") "359.31"
$ws.Range("E5").Value = "  +2.35%  "
Set-TextValue $ws.Range("D6") "109.47"
$ws.Range("E6").Value = "  -2.03%  "
Set-TextValue $ws.Range("D7") "0.565"
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("E8").Value = "  +0.01%  "
Set-TextValue $ws.Range("D9") "0.593"
$ws.Range("E9").Value = "  -0.69%  "
Set-TextValue $ws.Range("D10") "40.01"
$ws.Range("E10").Value = "  -3.00%  "
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("E13").Value = "  -1.51%  "
Set-TextValue $ws.Range("D14") "7.58"
$ws.Range("E14").Value = "  -1.86%  "
Set-TextValue $ws.Range("D15") "3.225.86"
$ws.Range("E15").Value = "  -0.68%  "
Set-TextValue $ws.Range("D16") "2.794.11"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("E17").Value = "  +8.02%  "
Set-TextValue $ws.Range("D18") "51.880.02"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("E19").Value = "  -0.58%  "
Set-TextValue $ws.Range("D20") "3.14"
$ws.Range("E20").Value = "  -1.14%  "
Set-TextValue $ws.Range("D21") "12.99"
$ws.Range("E21").Value = "  -2.24%  "
$ws.Range("E22").Value = "  -0.60%  "
Set-TextValue $ws.Range("D23") "274.48"
$ws.Range("E23").Value = "  +2.05%  "
Set-TextValue $ws.Range("D24") "70.26"
$ws.Range("E24").Value = "  +1.18%  "
Set-TextValue $ws.Range("D25") "2.75"
$ws.Range("E25").Value = "  +0.09%  "
Set-TextValue $ws.Range("D26") "26.71"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("E27").Value = "  +0.05%  "
Set-TextValue $ws.Range("D28") "2.28"
$ws.Range("E28").Value = "  +1.93%  "
Set-TextValue $ws.Range("D29") "10.19"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("E30").Value = "  +4.51%  "
$ws.Range("B31").Value = "OKB"
$ws.Range("C31").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D31") "51.66"
$ws.Range("E31").Value = "  +2.16%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D32") "34.62"
$ws.Range("E32").Value = "  +2.75%  "
$ws.Range("B33").Value = "VeChain"
$ws.Range("C33").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D33") "0.0462"
$ws.Range("E33").Value = "  +2.67%  "
Set-TextValue $ws.Range("D34") "5.72"
$ws.Range("E34").Value = "  -0.92%  "
Set-TextValue $ws.Range("D35") "0.0847"
$ws.Range("E35").Value = "  +3.56%  "
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("E39").Value = "  -2.41%  "
Set-TextValue $ws.Range("D40") "18.08"
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D43") "122.34"
$ws.Range("E43").Value = "  -2.85%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D44") "2.25"
$ws.Range("E44").Value = "  -1.72%  "
Set-TextValue $ws.Range("D45") "22.05"
$ws.Range("E45").Value = "  -7.28%  "
Set-TextValue $ws.Range("D46") "2.076.42"
$ws.Range("E46").Value = "  +0.38%  "
Set-TextValue $ws.Range("D47") "3.26"
$ws.Range("E47").Value = "  -1.46%  "
Set-TextValue $ws.Range("D48") "2.20"
$ws.Range("E48").Value = "  -3.54%  "
$ws.Range("E49").Value = "  +1.71%  "
Set-TextValue $ws.Range("D50") "0.935"
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("E51").Value = "  +0.72%  "
